# Fruta / hortaliza, semanal
# Inserts a new week's worth of data (2 rows) for "Betarraga" at
# Vega Monumental Concepción, right after the existing 2021-04-13 pair
# (rows 255-256), pushing the rest of the weekly pairs down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 257-258 (existing rows 257+ shift down by 2).
$ws.Rows.Item(257).Resize(2).Insert()

# Row 257: new "Primera" quality entry for the inserted week.
$ws.Range("A257").Value = 11
$ws.Range("B257").Value = "Vega Monumental Concepción"
$ws.Range("C257").Value = "Bíobío"
$ws.Range("D257").Value = 44610
$ws.Range("E257").Value = 8
$ws.Range("F257").Value = 100114014
$ws.Range("G257").Value = "Betarraga"
$ws.Range("H257").Value = "Sin especificar"
$ws.Range("I257").Value = "Primera"
$ws.Range("J257").Value = 800
$ws.Range("K257").Value = 600
$ws.Range("L257").Value = 700
$ws.Range("M257").Value = 650
$ws.Range("N257").Value = "`$/paquete 5 unidades"
$ws.Range("O257").Value = "Región Metropolitana"
$ws.Range("P257").Value = 130
$ws.Range("Q257").Value = 5
$ws.Range("R257").Value = "Hortaliza"

# Row 258: new "Segunda" quality entry for the inserted week.
$ws.Range("A258").Value = 11
$ws.Range("B258").Value = "Vega Monumental Concepción"
$ws.Range("C258").Value = "Bíobío"
$ws.Range("D258").Value = 44610
$ws.Range("E258").Value = 8
$ws.Range("F258").Value = 100114014
$ws.Range("G258").Value = "Betarraga"
$ws.Range("H258").Value = "Sin especificar"
$ws.Range("I258").Value = "Segunda"
$ws.Range("J258").Value = 400
$ws.Range("K258").Value = 500
$ws.Range("L258").Value = 500
$ws.Range("M258").Value = 500
$ws.Range("N258").Value = "`$/paquete 5 unidades"
$ws.Range("O258").Value = "Región Metropolitana"
$ws.Range("P258").Value = 100
$ws.Range("Q258").Value = 5
$ws.Range("R258").Value = "Hortaliza"

Write-Output "Inserted rows 257-258; sheet now has $($ws.Rows.Count) tracked rows"
